$d = $word.ActiveDocument

# --- 1. Heading paragraph "3.1. RESTRICCIONES DEL SOFTWARE" ---
# Bump the font size (sz/szCs = 32 half-points = 16pt) on the paragraph
# mark and on both runs that make up the heading.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Font.Size = 16
$p1.Range.Font.SizeBi = 16

# Wrap the whole heading paragraph with the "_GoBack" bookmark. Word only
# ever keeps a single "_GoBack" bookmark, so re-adding it here relocates
# it away from its old position at the end of the document.
$d.Bookmarks.Add("_GoBack", $d.Range($p1.Range.Start, $p1.Range.End))

# --- 2. Merge the split "e-mail" run back into a single run ---
# The two runs "...el mismo e-" and "mail registrado." were only split so a
# (now relocated) bookmark could sit between them; re-assert the combined
# text so the content collapses back into one run.
$d.Content.Find.Execute(
    "No puede existir mas de un usuario con el mismo e-mail registrado.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "No puede existir mas de un usuario con el mismo e-mail registrado.", 2
) | Out-Null
